$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - invoice date + QR code
$ws.Range("B5").Value = "02/06/2024"
$ws.Range("L5").Value = "QR02062024LP01"

# Row 10 - pharmacy / distributor name
$ws.Range("B10").Value = "LM PHARMACY"
$ws.Range("D10").Value = "KINARY"
$ws.Range("E10").Value = "LM PHARMACY"

# Row 11 - address line 1
$ws.Range("B11").Value = "900W SAM HOUSTON BLVD STE 3"
$ws.Range("D11").Value = "152.35 TENTH AVE"
$ws.Range("E11").Value = "900W SAM HOUSTON BLVD STE 3"

# Row 12 - address line 2 (city/state/zip)
$ws.Range("B12").Value = "PHARR, TX 78577"
$ws.Range("D12").Value = "WHITESTONE, NY,11357"
$ws.Range("E12").Value = "PHARR, TX 78577"

# Row 13 - account number
$ws.Range("D13").Value = "Account#: "

# Row 14 - phone numbers
$ws.Range("B14").Value = "Phone: 201-595-1234, fax: "
$ws.Range("D14").Value = "Phone: 718-767-1234/ 888-527-6806"
$ws.Range("E14").Value = "Phone: 201-595-1234, fax: "

# Row 15 - DEA numbers
$ws.Range("B15").Value = "DEA: BL5101009, Exp: 08/31/2023"
$ws.Range("D15").Value = "DEA: RK0416900"
$ws.Range("E15").Value = "DEA: BL5101009, Exp: 08/31/2023"

# Row 20 - line item 1
$ws.Range("B20").Value = "2315574603"
$ws.Range("C20").Value = "Avet Pharmaceuticals Inc."
$ws.Range("D20").Value = "Rasagiline mesylate"
$ws.Range("E20").Value = ".5 mg/1"
$ws.Range("G20").Value = "RCY01AD6"
$ws.Range("H20").Value = "12/23/31"
$ws.Range("I20").Value = "30 CT"

# Row 21 - line item 2
$ws.Range("B21").Value = "6275651818"
$ws.Range("C21").Value = "Sun Pharmaceutical Industries, Inc."
$ws.Range("D21").Value = "CARBIDOPA AND LEVODOPA"
$ws.Range("E21").Value = "25 mg/1"
$ws.Range("G21").Value = "HAD1849A"
$ws.Range("H21").Value = "04/24/30"
$ws.Range("I21").Value = "1000 CT"
